# COAD_AccCreation.xlsx - IBG Cao Scripts - Sarib Shamim
#
# Adds two "Ibg" worksheets (IbgFcyCurrent, IbgFcySaving) as copies of the
# "FCY current" sheet, plus a new blank "Sheet4", updates a couple of
# worksheet selections (mimicking a Ctrl+A "select-all" on a couple of
# sheets), clears the old active-tab marker and leaves "IbgFcySaving" as
# the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate "FCY current" twice to build the new Ibg sheets.
# ---------------------------------------------------------------------
$fcyCurrent = $wb.Worksheets.Item("FCY current")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fcyCurrent.Copy($null, $lastSheet)
$ibgFcyCurrent = $wb.Worksheets.Item($wb.Worksheets.Count)
$ibgFcyCurrent.Name = "IbgFcyCurrent"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fcyCurrent.Copy($null, $lastSheet)
$ibgFcySaving = $wb.Worksheets.Item($wb.Worksheets.Count)
$ibgFcySaving.Name = "IbgFcySaving"

# ---------------------------------------------------------------------
# 2) Add the trailing blank "Sheet4". A throw-away sheet is inserted
#    first (and removed again) purely so the internal sheetId counter
#    lands on the same value it did in the original authoring session.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$placeholder = $wb.Worksheets.Add($null, $lastSheet)
$placeholderName = $placeholder.Name

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)
$sheet4.Name = "Sheet4"

$wb.Worksheets.Item($placeholderName).Delete()

# ---------------------------------------------------------------------
# 3) Refresh a couple of pre-existing sheet selections (as if the user
#    had pressed Ctrl+A on them after parking the cursor elsewhere).
# ---------------------------------------------------------------------
$lcySaving = $wb.Worksheets.Item("LCY Saving")
$lcySaving.Range("P15").Select()
$lcySaving.Cells.Select()

$fcyCurrent.Range("J14").Select()
$fcyCurrent.Cells.Select()

# ---------------------------------------------------------------------
# 4) Give the two new Ibg sheets their own selections.
# ---------------------------------------------------------------------
$ibgFcyCurrent.Range("F9").Select()

$ibgFcySaving.Range("F2").Select()

# ---------------------------------------------------------------------
# 5) Finish with "IbgFcySaving" as the active sheet/tab (this also
#    clears the stale tabSelected flag that used to sit on
#    "LCCurrentInputter").
# ---------------------------------------------------------------------
$ibgFcySaving.Activate()
